# "Save plots for multi p" — add the 11th run sheet
# (Run_11_Chanells_10_BER_0.01) with its Channel/k/n/p/dmin/R/Data_Portion/
# Bitrate/"BER with ECC" table, appended after the existing ten run sheets.

$wb = $excel.ActiveWorkbook

# --- add the new worksheet as the last tab -------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Run_11_Chanells_10_BER_0.01"

# --- table data (header + 10 data rows, columns A:I) ----------------------
$rowsData = @(
    @("Channel", "k", "n", "p", "dmin", "R", "Data_Portion", "Bitrate", "BER with ECC"),
    @(1, 2, 3, 0.002, 2, 0.66666666666666663, 0.16055045871559631, 666666.66666666663, 0.0028994201159768048),
    @(2, 2, 5, 0.0073333333333333332, 3, 0.40000000000000002, 0.096330275229357804, 400000, 0.00019996000799840031),
    @(3, 2, 5, 0.012666666666666666, 3, 0.40000000000000002, 0.096330275229357804, 400000, 0.00099980003999200159),
    @(4, 2, 5, 0.018000000000000002, 3, 0.40000000000000002, 0.096330275229357804, 400000, 0.0010997800439912018),
    @(5, 2, 5, 0.023333333333333331, 3, 0.40000000000000002, 0.096330275229357804, 400000, 0.0035992801439712059),
    @(6, 2, 5, 0.028666666666666667, 3, 0.40000000000000002, 0.096330275229357804, 400000, 0.0043991201759648072),
    @(7, 2, 5, 0.034000000000000002, 3, 0.40000000000000002, 0.096330275229357804, 400000, 0.0054989002199560084),
    @(8, 2, 5, 0.039333333333333338, 3, 0.40000000000000002, 0.096330275229357804, 400000, 0.0072985402919416116),
    @(9, 2, 5, 0.044666666666666667, 3, 0.40000000000000002, 0.096330275229357804, 400000, 0.0099980003999200154),
    @(10, 2, 7, 0.050000000000000003, 4, 0.2857142857142857, 0.068807339449541274, 285714.28571428568, 0.009398120375924815)
)

$rowCount = $rowsData.Count
$colCount = $rowsData[0].Count

$values = New-Object 'object[,]' $rowCount, $colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $values[$r, $c] = $rowsData[$r][$c]
    }
}

$target = $newSheet.Range($newSheet.Cells.Item(1, 1), $newSheet.Cells.Item($rowCount, $colCount))
$target.Value = $values

# --- keep the first sheet as the active tab (unchanged from before) ------
$wb.Worksheets.Item(1).Activate()

# --- force a full recalculation on next open ------------------------------
$excel.CalculateFullRebuild()
